$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 436
$ws.Range("I5").Value = 424
$ws.Range("J5").Value = 460
$ws.Range("K5").Value = 424
$ws.Range("L5").Value = 460
$ws.Range("M5").Value = -309
$ws.Range("N5").Value = -690

$ws.Range("H11").Value = 167.5
$ws.Range("I11").Value = 167.5
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 167.5
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -27.5

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H33").Value = 213.47058
$ws.Range("I33").Value = 255
$ws.Range("J33").Value = 113.8
$ws.Range("K33").Value = 255
$ws.Range("L33").Value = 113.8
$ws.Range("M33").Value = -26

$ws.Range("H42").Value = 456.8125
$ws.Range("I42").Value = 126.111115
$ws.Range("J42").Value = 882
$ws.Range("K42").Value = 378.333345
$ws.Range("L42").Value = 2646
$ws.Range("M42").Value = -148.333345

$ws.Range("H93").Value = 99998
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 99998
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 99998
$ws.Range("N93").Value = -104990

$ws.Range("H106").Value = 2805.3333
$ws.Range("I106").Value = 2464
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 2464
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -1833

$ws.Range("H107").Value = 958.64703
$ws.Range("I107").Value = 958.64703
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 958.64703
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 961.35297
$ws.Range("N107").ClearContents()

$ws.Range("H110").Value = 52924.5
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 52924.5
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 52924.5
$ws.Range("N110").Value = -61104.5

$ws.Range("H129").Value = 142858960
$ws.Range("I129").Value = 200000750
$ws.Range("J129").Value = 4448.5
$ws.Range("K129").Value = 600002250
$ws.Range("L129").Value = 13345.5
$ws.Range("M129").Value = -599997250
$ws.Range("N129").Value = -23345.5

$ws.Range("H132").Value = 16964.666
$ws.Range("I132").Value = 5619.3335
$ws.Range("J132").Value = 51000.668
$ws.Range("K132").Value = 16858.0005
$ws.Range("L132").Value = 153002.004
$ws.Range("M132").Value = -14328.0005

$ws.Range("H133").Value = 108497.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 108497.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 108497.5
$ws.Range("N133").Value = -118617.5

$ws.Range("H137").Value = 2582.2083
$ws.Range("I137").Value = 2532.0476
$ws.Range("J137").Value = 2933.3333
$ws.Range("K137").Value = 7596.1428
$ws.Range("L137").Value = 8799.999899999999
$ws.Range("M137").Value = -5046.1428
$ws.Range("N137").Value = -13899.9999

$ws.Range("H141").Value = 3034.3333
$ws.Range("I141").Value = 1737.8572
$ws.Range("J141").Value = 5627.2856
$ws.Range("K141").Value = 5213.571599999999
$ws.Range("L141").Value = 16881.8568
$ws.Range("M141").Value = -33.57159999999931
$ws.Range("N141").Value = -27241.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16653.9
$ws.Range("I32").Value = 16653.9
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 16653.9
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -16366.9

$ws.Range("H43").Value = 16989.75
$ws.Range("I43").Value = 10171
$ws.Range("J43").Value = 19262.666
$ws.Range("K43").Value = 10171
$ws.Range("L43").Value = 19262.666
$ws.Range("M43").Value = -9858
$ws.Range("N43").Value = -19888.666

$ws.Range("H61").Value = 5559
$ws.Range("I61").Value = 2650.6667
$ws.Range("J61").Value = 12270.538
$ws.Range("K61").Value = 2650.6667
$ws.Range("L61").Value = 12270.538
$ws.Range("M61").Value = -2438.6667

$ws.Range("H74").Value = 1859.4736
$ws.Range("I74").Value = 1859.4736
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1859.4736
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -985.4736

$ws.Range("H77").Value = 1859.4736
$ws.Range("I77").Value = 1859.4736
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 9297.368
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -4929.368

$ws.Range("H88").Value = 1640.45
$ws.Range("I88").Value = 1372.1428
$ws.Range("J88").Value = 1784.9231
$ws.Range("K88").Value = 1372.1428
$ws.Range("L88").Value = 1784.9231
$ws.Range("M88").Value = -966.1428000000001
$ws.Range("N88").Value = -2596.9231

$ws.Range("H91").Value = 1640.45
$ws.Range("I91").Value = 1372.1428
$ws.Range("J91").Value = 1784.9231
$ws.Range("K91").Value = 1372.1428
$ws.Range("L91").Value = 1784.9231
$ws.Range("M91").Value = 31.85719999999992
$ws.Range("N91").Value = -4592.9231

$ws.Range("H102").Value = 2367.25
$ws.Range("I102").Value = 1465
$ws.Range("J102").Value = 3630.4
$ws.Range("K102").Value = 1465
$ws.Range("L102").Value = 3630.4
$ws.Range("M102").Value = 157

$ws.Range("H110").Value = 1963.1538
$ws.Range("I110").Value = 1091.125
$ws.Range("J110").Value = 3358.4
$ws.Range("K110").Value = 1091.125
$ws.Range("L110").Value = 3358.4
$ws.Range("M110").Value = 953.875

$ws.Range("H122").Value = 2293.543
$ws.Range("I122").Value = 1946.125
$ws.Range("J122").Value = 5999.3335
$ws.Range("K122").Value = 5838.375
$ws.Range("L122").Value = 17998.0005
$ws.Range("M122").Value = -3388.375

$ws.Range("H132").Value = 1598.6531
$ws.Range("I132").Value = 1276.909
$ws.Range("J132").Value = 4430
$ws.Range("K132").Value = 3830.727
$ws.Range("L132").Value = 13290
$ws.Range("M132").Value = -1300.727

$ws.Range("H136").Value = 5559
$ws.Range("I136").Value = 2650.6667
$ws.Range("J136").Value = 12270.538
$ws.Range("K136").Value = 7952.000100000001
$ws.Range("L136").Value = 36811.614
$ws.Range("M136").Value = -5402.000100000001

$ws.Range("H139").Value = 104446.664
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 104446.664
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 104446.664
$ws.Range("N139").Value = -114726.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 28332.666
$ws.Range("I26").Value = 28332.666
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 28332.666
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -28040.666

$ws.Range("H86").Value = 6038.4443
$ws.Range("I86").Value = 9612.429
$ws.Range("J86").Value = 2189.5386
$ws.Range("K86").Value = 9612.429
$ws.Range("L86").Value = 2189.5386
$ws.Range("M86").Value = -8489.429
$ws.Range("N86").Value = -4435.5386

$ws.Range("H89").Value = 6038.4443
$ws.Range("I89").Value = 9612.429
$ws.Range("J89").Value = 2189.5386
$ws.Range("K89").Value = 48062.145
$ws.Range("L89").Value = 10947.693
$ws.Range("M89").Value = -42446.145
$ws.Range("N89").Value = -22179.693

$ws.Range("H105").Value = 3360.0286
$ws.Range("I105").Value = 2023.7931
$ws.Range("J105").Value = 9818.5
$ws.Range("K105").Value = 2023.7931
$ws.Range("L105").Value = 9818.5
$ws.Range("M105").Value = -276.7931000000001
$ws.Range("N105").Value = -13312.5

$ws.Range("H130").Value = 74999.5
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 74999.5
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 74999.5
$ws.Range("N130").Value = -85039.5

$ws.Range("H134").Value = 4564.6606
$ws.Range("I134").Value = 1832.1471
$ws.Range("J134").Value = 8787.637000000001
$ws.Range("K134").Value = 5496.4413
$ws.Range("L134").Value = 26362.911
$ws.Range("M134").Value = -2961.4413

$ws.Range("H135").Value = 90486.664
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 90486.664
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 90486.664
$ws.Range("N135").Value = -100626.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 814
$ws.Range("I2").Value = 814
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 814
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -701

$ws.Range("H4").Value = 500000500
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 1000000000
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 1000000000
$ws.Range("M4").Value = -888
$ws.Range("N4").Value = -1000000224

$ws.Range("H16").Value = 1656.7059
$ws.Range("I16").Value = 1351.9166
$ws.Range("J16").Value = 2388.2
$ws.Range("K16").Value = 1351.9166
$ws.Range("L16").Value = 2388.2
$ws.Range("M16").Value = -1064.9166

$ws.Range("H22").Value = 529.5
$ws.Range("I22").Value = 492
$ws.Range("J22").Value = 617
$ws.Range("K22").Value = 492
$ws.Range("L22").Value = 617
$ws.Range("M22").Value = -142

$ws.Range("H31").Value = 2565.037
$ws.Range("I31").Value = 1398.0454
$ws.Range("J31").Value = 7699.8
$ws.Range("K31").Value = 1398.0454
$ws.Range("L31").Value = 7699.8
$ws.Range("M31").Value = -1103.0454

$ws.Range("H34").Value = 2565.037
$ws.Range("I34").Value = 1398.0454
$ws.Range("J34").Value = 7699.8
$ws.Range("K34").Value = 1398.0454
$ws.Range("L34").Value = 7699.8
$ws.Range("M34").Value = -1196.0454

$ws.Range("H48").Value = 11821.8
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 11821.8
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 11821.8
$ws.Range("N48").Value = -12773.8

$ws.Range("H68").Value = 60295
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 60295
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 60295
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -61793

$ws.Range("H71").Value = 60295
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 60295
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 180885
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -188373

$ws.Range("H99").Value = 2788.5557
$ws.Range("I99").Value = 2585.7144
$ws.Range("J99").Value = 3498.5
$ws.Range("K99").Value = 2585.7144
$ws.Range("L99").Value = 3498.5
$ws.Range("M99").Value = -1087.7144
$ws.Range("N99").Value = -6494.5

$ws.Range("H113").Value = 1656.7059
$ws.Range("I113").Value = 1351.9166
$ws.Range("J113").Value = 2388.2
$ws.Range("K113").Value = 1351.9166
$ws.Range("L113").Value = 2388.2
$ws.Range("M113").Value = 818.0834

$ws.Range("H126").Value = 2788.5557
$ws.Range("I126").Value = 2585.7144
$ws.Range("J126").Value = 3498.5
$ws.Range("K126").Value = 7757.1432
$ws.Range("L126").Value = 10495.5
$ws.Range("M126").Value = -5287.1432
$ws.Range("N126").Value = -15435.5

$ws.Range("H132").Value = 4004510.8
$ws.Range("I132").Value = 6672034.5
$ws.Range("J132").Value = 3225
$ws.Range("K132").Value = 20016103.5
$ws.Range("L132").Value = 9675
$ws.Range("M132").Value = -20013573.5
$ws.Range("N132").Value = -14735

$ws.Range("H134").Value = 4785.278
$ws.Range("I134").Value = 1891.875
$ws.Range("J134").Value = 7100
$ws.Range("K134").Value = 5675.625
$ws.Range("L134").Value = 21300
$ws.Range("M134").Value = -3140.625

$ws.Range("H141").Value = 95480.75
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 95480.75
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 95480.75
$ws.Range("N141").Value = -105840.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12022015
$ws.Range("I4").Value = 11554620
$ws.Range("J4").Value = 14242141
$ws.Range("K4").Value = 34663860
$ws.Range("L4").Value = 42726423
$ws.Range("M4").Value = -34663748
$ws.Range("N4").Value = -42726647

$ws.Range("H33").Value = 363.3846
$ws.Range("I33").Value = 207.4
$ws.Range("J33").Value = 460.875
$ws.Range("K33").Value = 1244.4
$ws.Range("L33").Value = 2765.25
$ws.Range("M33").Value = -961.4000000000001
$ws.Range("N33").Value = -3331.25

$ws.Range("H122").Value = 645.6818
$ws.Range("I122").Value = 315.66666
$ws.Range("J122").Value = 769.4375
$ws.Range("K122").Value = 2840.99994
$ws.Range("L122").Value = 6924.9375
$ws.Range("M122").Value = -390.9999399999997
$ws.Range("N122").Value = -11824.9375

$ws.Range("H131").Value = 1654.8695
$ws.Range("I131").Value = 1031.2
$ws.Range("J131").Value = 1828.1111
$ws.Range("K131").Value = 3093.6
$ws.Range("L131").Value = 5484.3333
$ws.Range("M131").Value = 1946.4
$ws.Range("N131").Value = -15564.3333

$ws.Range("H137").Value = 2516.611
$ws.Range("I137").Value = 1007
$ws.Range("J137").Value = 3477.2727
$ws.Range("K137").Value = 3021
$ws.Range("L137").Value = 10431.8181
$ws.Range("M137").Value = 2079
$ws.Range("N137").Value = -20631.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 224.86667
$ws.Range("I2").Value = 183.5
$ws.Range("J2").Value = 272.14285
$ws.Range("K2").Value = 183.5
$ws.Range("L2").Value = 272.14285
$ws.Range("M2").Value = -70.5

$ws.Range("H70").Value = 4315.8667
$ws.Range("I70").Value = 3994.3333
$ws.Range("J70").Value = 4798.1665
$ws.Range("K70").Value = 3994.3333
$ws.Range("L70").Value = 4798.1665
$ws.Range("M70").Value = -3724.3333
$ws.Range("N70").Value = -5338.1665

$ws.Range("H73").Value = 4315.8667
$ws.Range("I73").Value = 3994.3333
$ws.Range("J73").Value = 4798.1665
$ws.Range("K73").Value = 3994.3333
$ws.Range("L73").Value = 4798.1665
$ws.Range("M73").Value = -3058.3333
$ws.Range("N73").Value = -6670.1665

$ws.Range("H103").Value = 38000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 38000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 38000
$ws.Range("N103").Value = -40344

$ws.Range("H112").Value = 84997.5
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 84997.5
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 84997.5
$ws.Range("N112").Value = -87213.5

$ws.Range("H113").Value = 1682
$ws.Range("I113").Value = 1682
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1682
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 488
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 18522062
$ws.Range("I132").Value = 18522062
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 55566186
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -55563656
$ws.Range("N132").ClearContents()

$ws.Range("H140").Value = 95000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 95000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 95000
$ws.Range("N140").Value = -105360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -888

$ws.Range("H35").Value = 2343.2727
$ws.Range("I35").Value = 2343.2727
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2343.2727
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -2007.2727
$ws.Range("N35").ClearContents()

$ws.Range("H36").Value = 74857.5
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 74857.5
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 74857.5
$ws.Range("N36").Value = -75981.5

$ws.Range("H46").Value = 10999.083
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 10999.083
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 10999.083
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -11375.083

$ws.Range("H55").Value = 186.28
$ws.Range("I55").Value = 116.61539
$ws.Range("J55").Value = 261.75
$ws.Range("K55").Value = 116.61539
$ws.Range("L55").Value = 261.75
$ws.Range("M55").Value = 56.38461
$ws.Range("N55").Value = -607.75

$ws.Range("H132").Value = 2898.255
$ws.Range("I132").Value = 2847.75
$ws.Range("J132").Value = 3081.9092
$ws.Range("K132").Value = 8543.25
$ws.Range("L132").Value = 9245.7276
$ws.Range("M132").Value = -6013.25

$ws.Range("H136").Value = 4013.2273
$ws.Range("I136").Value = 3999.9412
$ws.Range("J136").Value = 4058.4
$ws.Range("K136").Value = 11999.8236
$ws.Range("L136").Value = 12175.2
$ws.Range("M136").Value = -9449.8236

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 34994.5
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 34994.5
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 34994.5
$ws.Range("N47").Value = -36138.5

$ws.Range("H81").Value = 64957.824
$ws.Range("I81").Value = 114723.445
$ws.Range("J81").Value = 8971.5
$ws.Range("K81").Value = 229446.89
$ws.Range("L81").Value = 17943
$ws.Range("M81").Value = -228385.89
$ws.Range("N81").Value = -20065

$ws.Range("H84").Value = 64957.824
$ws.Range("I84").Value = 114723.445
$ws.Range("J84").Value = 8971.5
$ws.Range("K84").Value = 1147234.45
$ws.Range("L84").Value = 89715
$ws.Range("M84").Value = -1141930.45
$ws.Range("N84").Value = -100323

$ws.Range("H103").Value = 43666.668
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 43666.668
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 43666.668
$ws.Range("N103").Value = -46010.668

$ws.Range("H123").Value = 114989
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 114989
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 114989
$ws.Range("N123").Value = -124789

$ws.Range("H132").Value = 2772.3572
$ws.Range("I132").Value = 2876.3914
$ws.Range("J132").Value = 2293.8
$ws.Range("K132").Value = 8629.174199999999
$ws.Range("L132").Value = 6881.400000000001
$ws.Range("M132").Value = -6099.174199999999
$ws.Range("N132").Value = -11941.4

$ws.Range("H136").Value = 1510.0555
$ws.Range("I136").Value = 1684.8462
$ws.Range("J136").Value = 1055.6
$ws.Range("K136").Value = 5054.5386
$ws.Range("L136").Value = 3166.8
$ws.Range("M136").Value = -2504.5386
